$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.443.65"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "3.715.65"
$ws.Range("E3").Value = "  -0.45%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "611.79"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.18"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.67%  "
$ws.Range("D7").Value = "3.712.43"
$ws.Range("E7").Value = "  -0.52%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.528"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.164"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.33%  "
$ws.Range("E11").Value = "  +2.56%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.478"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -4.83%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "39.60"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -3.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000252"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.57%  "
$ws.Range("D15").Value = "4.330.86"
$ws.Range("E15").Value = "  -0.67%  "
$ws.Range("D16").Value = "3.714.98"
$ws.Range("E16").Value = "  -0.63%  "
$ws.Range("D17").Value = "69.494.30"
$ws.Range("E17").Value = "  -0.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.120"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -2.63%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.46"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "499.94"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -3.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.24"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.11"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.716"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.58"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +3.60%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.92"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.21"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.86"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -5.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000133"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +6.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.996"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.39%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.45"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.53%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.02"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.31%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.90"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.25"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.96%  "
$ws.Range("E34").Value = "  -2.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.05"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.98%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.07"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.347"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.66%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.138"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +4.58%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.02"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +11.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.05"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -6.32%  "
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "49.66"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.23%  "
$ws.Range("B43").Value = "Arweave"
$ws.Range("C43").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "45.47"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "433.20"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.53"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.70%  "
$ws.Range("D46").Value = "2.940.36"
$ws.Range("E46").Value = "  -3.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0360"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "138.78"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.64%  "
$ws.Range("B49").Value = "USDe"
$ws.Range("C49").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "27.01"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.45"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.14%  "
